## "qualifier 1 completed with predictions"
##
## Row 66 on Sheet1 is the "Qualifier 1" (MI vs DC) row. Enter the six
## predicted scores (columns E, H, K, N, Q, T - one per participant's
## prediction) so the VLOOKUP/RANK formulas in D, G, J, M, P, S resolve to
## real point values instead of the empty "" placeholder.
##
## Also record this match in the "Most Wins" / prize tracker table
## (B81:G86 / L81:R86): column C gets each participant's updated win
## count and column G becomes a running total formula instead of a
## hand-typed number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Qualifier 1 (MI vs DC) predictions entered for each participant.
$ws.Range("E66").Value = 0
$ws.Range("H66").Value = 80
$ws.Range("K66").Value = 100
$ws.Range("N66").Value = 20
$ws.Range("Q66").Value = 40
$ws.Range("T66").Value = 60

# "Most wins" counts, updated now that qualifier 1 results are in.
$ws.Range("C81").Value = 0
$ws.Range("C82").Value = 0
$ws.Range("C83").Value = 10
$ws.Range("C84").Value = 7
$ws.Range("C85").Value = 3
$ws.Range("C86").Value = 3

# Column G becomes a live total (=SUM(C:F)) instead of a static number.
$ws.Range("G81").Formula = "=SUM(C81:F81)"
$ws.Range("G82:G86").Formula = "=SUM(C82:F82)"

$wb.Application.Calculate()
